$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04487447322743776
$ws.Range("D2").Value = 0.6671606830458927
$ws.Range("E2").Value = 0.1436091338164971
$ws.Range("F2").Value = 11.84800143897968
$ws.Range("G2").Value = 0.002756231402716107
$ws.Range("J2").Value = 0.06072616251083929
$ws.Range("M2").Value = 5.707467284634873
$ws.Range("N2").Value = 1.232680290454681
$ws.Range("C3").Value = 0.03993442494882515
$ws.Range("D3").Value = 0.6065466171190792
$ws.Range("E3").Value = 0.1252159035566152
$ws.Range("F3").Value = 11.53170165884114
$ws.Range("G3").Value = 0.002779715310771446
$ws.Range("J3").Value = 0.05697754077701944
$ws.Range("M3").Value = 5.373027205438262
$ws.Range("N3").Value = 1.20805798644318
$ws.Range("C4").Value = 0.03695355994240401
$ws.Range("D4").Value = 0.5701880976156417
$ws.Range("E4").Value = 0.1140304883569456
$ws.Range("F4").Value = 11.35143693870538
$ws.Range("G4").Value = 0.002794751389554391
$ws.Range("J4").Value = 0.05467329212353889
$ws.Range("M4").Value = 5.172333241996171
$ws.Range("N4").Value = 1.193599926557965
$ws.Range("C5").Value = 0.03575094595377948
$ws.Range("D5").Value = 0.5555732440484462
$ws.Range("E5").Value = 0.1094965091582338
$ws.Range("F5").Value = 11.28136076970674
$ws.Range("G5").Value = 0.002801035615004811
$ws.Range("J5").Value = 0.05373341244828822
$ws.Range("M5").Value = 5.091667530547312
$ws.Range("N5").Value = 1.18787577143128
$ws.Range("C6").Value = 0.03555195462810445
$ws.Range("D6").Value = 0.5531582228263687
$ws.Range("E6").Value = 0.1087450196825728
$ws.Range("F6").Value = 11.26992551570146
$ws.Range("G6").Value = 0.002802088633558214
$ws.Range("J6").Value = 0.05357728652268179
$ws.Range("M6").Value = 5.078339163721239
$ws.Range("N6").Value = 1.186935452702215
$ws.Range("C7").Value = 0.0369372933779033
$ws.Range("D7").Value = 0.5699901991146703
$ws.Range("E7").Value = 0.1139692475512391
$ws.Range("F7").Value = 11.35047832733522
$ws.Range("G7").Value = 0.002794835503327286
$ws.Range("J7").Value = 0.05466062043507947
$ws.Range("M7").Value = 5.171240892474742
$ws.Range("N7").Value = 1.193522047701379
$ws.Range("C8").Value = 0.04315963195382722
$ws.Range("D8").Value = 0.6460738188985147
$ws.Range("E8").Value = 0.1372428462979016
$ws.Range("F8").Value = 11.73597089140668
$ws.Range("G8").Value = 0.002764201718821886
$ws.Range("J8").Value = 0.05943402889132443
$ws.Range("M8").Value = 5.591154565858062
$ws.Range("N8").Value = 1.224054453831229
$ws.Range("C9").Value = 0.05582733406190243
$ws.Range("D9").Value = 0.8027784826419975
$ws.Range("E9").Value = 0.1838854520261748
$ws.Range("F9").Value = 12.60847610859201
$ws.Range("G9").Value = 0.002708939334839647
$ws.Range("J9").Value = 0.06878574483158673
$ws.Range("M9").Value = 6.454031337204839
$ws.Range("N9").Value = 1.289102379110204
$ws.Range("C10").Value = 0.06549296815953198
$ws.Range("D10").Value = 0.9235096141553072
$ws.Range("E10").Value = 0.2189820882321811
$ws.Range("F10").Value = 13.32932256742708
$ws.Range("G10").Value = 0.002671148105109849
$ws.Range("J10").Value = 0.07566880053773417
$ws.Range("M10").Value = 7.115783242430382
$ws.Range("N10").Value = 1.33997233598086
$ws.Range("C11").Value = 0.06998573192032609
$ws.Range("D11").Value = 0.9798926306015119
$ws.Range("E11").Value = 0.2351798834179561
$ws.Range("F11").Value = 13.67666656973466
$ws.Range("G11").Value = 0.002654537812737313
$ws.Range("J11").Value = 0.07880706776586521
$ws.Range("M11").Value = 7.423766432514185
$ws.Range("N11").Value = 1.363768810038039
$ws.Range("C12").Value = 0.07170237089847831
$ws.Range("D12").Value = 1.001475157823336
$ws.Range("E12").Value = 0.241351503253
$ws.Range("F12").Value = 13.81117510461598
$ws.Range("G12").Value = 0.00264832906800697
$ws.Range("J12").Value = 0.07999684809389151
$ws.Range("M12").Value = 7.541470293090981
$ws.Range("N12").Value = 1.372872934671619
$ws.Range("C13").Value = 0.07133195450819585
$ws.Range("D13").Value = 0.9968163232978782
$ws.Range("E13").Value = 0.2400205755916573
$ws.Range("F13").Value = 13.78207083134805
$ws.Range("G13").Value = 0.002649662659235321
$ws.Range("J13").Value = 0.07974053987162222
$ws.Range("M13").Value = 7.516071442415864
$ws.Range("N13").Value = 1.370908086765638
$ws.Range("C14").Value = 0.0701266444173001
$ws.Range("D14").Value = 0.9816634711872325
$ws.Range("E14").Value = 0.2356868405887127
$ws.Range("F14").Value = 13.68767191106048
$ws.Range("G14").Value = 0.002654025401558271
$ws.Range("J14").Value = 0.07890492174499997
$ws.Range("M14").Value = 7.433427941434161
$ws.Range("N14").Value = 1.364515957132994
$ws.Range("C15").Value = 0.06939040096899873
$ws.Range("D15").Value = 0.9724127174380328
$ws.Range("E15").Value = 0.2330373725264678
$ws.Range("F15").Value = 13.63024318360158
$ws.Range("G15").Value = 0.002656708211874873
$ws.Range("J15").Value = 0.07839327331894452
$ws.Range("M15").Value = 7.382949123077822
$ws.Range("N15").Value = 1.360612654866543
$ws.Range("C16").Value = 0.06520141353608722
$ws.Range("D16").Value = 0.9198560540665994
$ws.Range("E16").Value = 0.2179285718599857
$ws.Range("F16").Value = 13.30702912806458
$ws.Range("G16").Value = 0.002672245122234901
$ws.Range("J16").Value = 0.07546388110106506
$ws.Range("M16").Value = 7.095802191108135
$ws.Range("N16").Value = 1.338430254344672
$ws.Range("C17").Value = 0.06265721543593372
$ws.Range("D17").Value = 0.8880033460440586
$ws.Range("E17").Value = 0.2087223646900611
$ws.Range("F17").Value = 13.11385005535095
$ws.Range("G17").Value = 0.002681923719314721
$ws.Range("J17").Value = 0.07366886816010521
$ws.Range("M17").Value = 6.921480412520424
$ws.Range("N17").Value = 1.324988980425047
$ws.Range("C18").Value = 0.06120278007507807
$ws.Range("D18").Value = 0.8698187452525872
$ws.Range("E18").Value = 0.2034487633413136
$ws.Range("F18").Value = 13.00455545817476
$ws.Range("G18").Value = 0.002687545438585424
$ws.Range("J18").Value = 0.07263707205970604
$ws.Range("M18").Value = 6.821864202135941
$ws.Range("N18").Value = 1.317319687895605
$ws.Range("C19").Value = 0.06071181982183305
$ws.Range("D19").Value = 0.8636845250489387
$ws.Range("E19").Value = 0.2016667849476477
$ws.Range("F19").Value = 12.96785685835027
$ws.Range("G19").Value = 0.002689458346081335
$ws.Range("J19").Value = 0.07228782436387604
$ws.Range("M19").Value = 6.788245055686616
$ws.Range("N19").Value = 1.314733648565152
$ws.Range("C20").Value = 0.0629271155912221
$ws.Range("D20").Value = 0.8913798718412522
$ws.Range("E20").Value = 0.2097001177418178
$ws.Range("F20").Value = 13.13422496804446
$ws.Range("G20").Value = 0.002680887756532576
$ws.Range("J20").Value = 0.07385988077579952
$ws.Range("M20").Value = 6.939969540118881
$ws.Range("N20").Value = 1.326413443637222
$ws.Range("C21").Value = 0.07048024436950584
$ws.Range("D21").Value = 0.9861077711686335
$ws.Range("E21").Value = 0.2369586986234538
$ws.Range("F21").Value = 13.71531684594362
$ws.Range("G21").Value = 0.002652741773578181
$ws.Range("J21").Value = 0.07915032238555852
$ws.Range("M21").Value = 7.457672452392899
$ws.Range("N21").Value = 1.36639096921553
$ws.Range("C22").Value = 0.07550678456192372
$ws.Range("D22").Value = 1.049378274452863
$ws.Range("E22").Value = 0.2549967850891051
$ws.Range("F22").Value = 14.11254005034476
$ws.Range("G22").Value = 0.002634818931790406
$ws.Range("J22").Value = 0.08261620353128762
$ws.Range("M22").Value = 7.802341065549797
$ws.Range("N22").Value = 1.393059600742987
$ws.Range("C23").Value = 0.07281523819155211
$ws.Range("D23").Value = 1.015477715447048
$ws.Range("E23").Value = 0.2453475463475314
$ws.Range("F23").Value = 13.89887429567835
$ws.Range("G23").Value = 0.002644342296155148
$ws.Range("J23").Value = 0.08076551348156613
$ws.Range("M23").Value = 7.61777936015676
$ws.Range("N23").Value = 1.37877697542018
$ws.Range("C24").Value = 0.06280506812470321
$ws.Range("D24").Value = 0.8898529475997634
$ws.Range("E24").Value = 0.2092580162127149
$ws.Range("F24").Value = 13.12500798107999
$ws.Range("G24").Value = 0.002681355936427321
$ws.Range("J24").Value = 0.07377352343614518
$ws.Range("M24").Value = 6.931608730967298
$ws.Range("N24").Value = 1.325769262203409
$ws.Range("C25").Value = 0.05234308385320219
$ws.Range("D25").Value = 0.759480242876748
$ws.Range("E25").Value = 0.1711384345417457
$ws.Range("F25").Value = 12.35909055883485
$ws.Range("G25").Value = 0.002723386617549511
$ws.Range("J25").Value = 0.06625505199635029
$ws.Range("M25").Value = 6.216005967925952
$ws.Range("N25").Value = 1.270959272791515
